$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value2 = -0.6954620276105178
$ws.Range("D2").Value2 = -0.1880822818662423
$ws.Range("E2").Value2 = 0.02914056319740949
$ws.Range("F2").Value2 = 0.03420913027474753
$ws.Range("G2").Value2 = 0.05668549591695958
$ws.Range("H2").Value2 = 0.08800201154955861
$ws.Range("J2").Value2 = 0.1026585402629121
$ws.Range("K2").Value2 = -0.09569731657107264
$ws.Range("L2").Value2 = -0.15961449175092
$ws.Range("M2").Value2 = 0.02671256140689233
$ws.Range("N2").Value2 = 0.1245000835920851
$ws.Range("O2").Value2 = 0.02613467962342986
$ws.Range("P2").Value2 = 0.1228681942988932
$ws.Range("Q2").Value2 = -0.033236964184887
$ws.Range("R2").Value2 = 0.001600103591062275
$ws.Range("B3").Value2 = -0.6954620276105178
$ws.Range("D3").Value2 = 0.2748306532819472
$ws.Range("E3").Value2 = -0.04857978195050124
$ws.Range("F3").Value2 = 0.01255724259136151
$ws.Range("G3").Value2 = -0.01060018019543008
$ws.Range("H3").Value2 = -0.1003275785479076
$ws.Range("J3").Value2 = -0.09495840075131633
$ws.Range("K3").Value2 = 0.1508618296685586
$ws.Range("L3").Value2 = 0.09568134073151335
$ws.Range("M3").Value2 = -0.08960471633066344
$ws.Range("N3").Value2 = -0.1001432587328002
$ws.Range("O3").Value2 = 0.01048147561741413
$ws.Range("P3").Value2 = -0.05947472534720946
$ws.Range("Q3").Value2 = 0.07520151954958959
$ws.Range("R3").Value2 = 0.06304668397542296
$ws.Range("B4").Value2 = -0.1880822818662423
$ws.Range("C4").Value2 = 0.2748306532819472
$ws.Range("E4").Value2 = -0.5123597682070902
$ws.Range("F4").Value2 = 0.5658061753654416
$ws.Range("G4").Value2 = 0.1739951819953533
$ws.Range("H4").Value2 = -0.2945441680967816
$ws.Range("J4").Value2 = -0.2905985076269564
$ws.Range("K4").Value2 = 0.3122181758795772
$ws.Range("L4").Value2 = -0.2754181976148924
$ws.Range("M4").Value2 = -0.03784677341782061
$ws.Range("N4").Value2 = -0.32874929123365
$ws.Range("O4").Value2 = 0.06019978600222384
$ws.Range("P4").Value2 = -0.2782325710945733
$ws.Range("Q4").Value2 = 0.2684339346906659
$ws.Range("R4").Value2 = -0.2432149506306212
$ws.Range("B5").Value2 = 0.02914056319740949
$ws.Range("C5").Value2 = -0.04857978195050124
$ws.Range("D5").Value2 = -0.5123597682070902
$ws.Range("F5").Value2 = -0.8861524251917112
$ws.Range("G5").Value2 = -0.3860168196198112
$ws.Range("H5").Value2 = 0.488511382759929
$ws.Range("J5").Value2 = 0.5327311332647087
$ws.Range("K5").Value2 = -0.6725077643418644
$ws.Range("L5").Value2 = 0.6935453139197432
$ws.Range("M5").Value2 = -0.09650071735825808
$ws.Range("N5").Value2 = 0.5236914481283271
$ws.Range("O5").Value2 = -0.2827997122318284
$ws.Range("P5").Value2 = 0.4559128407195984
$ws.Range("Q5").Value2 = -0.6020672565082923
$ws.Range("R5").Value2 = 0.5893034284480509
$ws.Range("B6").Value2 = 0.03420913027474753
$ws.Range("C6").Value2 = 0.01255724259136151
$ws.Range("D6").Value2 = 0.5658061753654416
$ws.Range("E6").Value2 = -0.8861524251917112
$ws.Range("G6").Value2 = 0.3335582946014879
$ws.Range("H6").Value2 = -0.5229843227986096
$ws.Range("J6").Value2 = -0.5422989646417723
$ws.Range("K6").Value2 = 0.4056758208164513
$ws.Range("L6").Value2 = -0.6304082599202043
$ws.Range("M6").Value2 = -0.001045408961199361
$ws.Range("N6").Value2 = -0.5212886093024264
$ws.Range("O6").Value2 = 0.270433041677948
$ws.Range("P6").Value2 = -0.446859052656183
$ws.Range("Q6").Value2 = 0.3654188162483002
$ws.Range("R6").Value2 = -0.4870047038530387
$ws.Range("B7").Value2 = 0.05668549591695958
$ws.Range("C7").Value2 = -0.01060018019543008
$ws.Range("D7").Value2 = 0.1739951819953533
$ws.Range("E7").Value2 = -0.3860168196198112
$ws.Range("F7").Value2 = 0.3335582946014879
$ws.Range("H7").Value2 = 0.2282144331985945
$ws.Range("J7").Value2 = 0.2017522194473726
$ws.Range("K7").Value2 = 0.3703485431360334
$ws.Range("L7").Value2 = -0.07823792394983131
$ws.Range("M7").Value2 = 0.813853831783909
$ws.Range("N7").Value2 = 0.1591991380261942
$ws.Range("O7").Value2 = -0.1978194138819634
$ws.Range("P7").Value2 = 0.1469122152520484
$ws.Range("Q7").Value2 = 0.3630575348992239
$ws.Range("R7").Value2 = -0.07265664635526758
$ws.Range("B8").Value2 = 0.08800201154955861
$ws.Range("C8").Value2 = -0.1003275785479076
$ws.Range("D8").Value2 = -0.2945441680967816
$ws.Range("E8").Value2 = 0.488511382759929
$ws.Range("F8").Value2 = -0.5229843227986096
$ws.Range("G8").Value2 = 0.2282144331985945
$ws.Range("J8").Value2 = 0.9792641151560048
$ws.Range("K8").Value2 = -0.05211753778901478
$ws.Range("L8").Value2 = 0.6829028338744026
$ws.Range("M8").Value2 = 0.48402467427776
$ws.Range("N8").Value2 = 0.7269516643205035
$ws.Range("O8").Value2 = -0.2839099923732194
$ws.Range("P8").Value2 = 0.6208983952960462
$ws.Range("Q8").Value2 = 0.04958524081280442
$ws.Range("R8").Value2 = 0.527323953863141
$ws.Range("B10").Value2 = 0.1026585402629121
$ws.Range("C10").Value2 = -0.09495840075131633
$ws.Range("D10").Value2 = -0.2905985076269564
$ws.Range("E10").Value2 = 0.5327311332647087
$ws.Range("F10").Value2 = -0.5422989646417723
$ws.Range("G10").Value2 = 0.2017522194473726
$ws.Range("H10").Value2 = 0.9792641151560048
$ws.Range("K10").Value2 = -0.09840861227319085
$ws.Range("L10").Value2 = 0.7164471536491567
$ws.Range("M10").Value2 = 0.4649076623324462
$ws.Range("N10").Value2 = 0.7604571890318254
$ws.Range("O10").Value2 = -0.2806547402725821
$ws.Range("P10").Value2 = 0.6602733475555668
$ws.Range("Q10").Value2 = 0.007045408321917977
$ws.Range("R10").Value2 = 0.5855897097062018
$ws.Range("B11").Value2 = -0.09569731657107264
$ws.Range("C11").Value2 = 0.1508618296685586
$ws.Range("D11").Value2 = 0.3122181758795772
$ws.Range("E11").Value2 = -0.6725077643418644
$ws.Range("F11").Value2 = 0.4056758208164513
$ws.Range("G11").Value2 = 0.3703485431360334
$ws.Range("H11").Value2 = -0.05211753778901478
$ws.Range("J11").Value2 = -0.09840861227319085
$ws.Range("L11").Value2 = -0.2857460969514533
$ws.Range("M11").Value2 = 0.3130512932906556
$ws.Range("N11").Value2 = -0.06611384108455527
$ws.Range("O11").Value2 = 0.03376595019548962
$ws.Range("P11").Value2 = -0.05246523971468067
$ws.Range("Q11").Value2 = 0.8799916556406883
$ws.Range("R11").Value2 = -0.2643290939009635
$ws.Range("B12").Value2 = -0.15961449175092
$ws.Range("C12").Value2 = 0.09568134073151335
$ws.Range("D12").Value2 = -0.2754181976148924
$ws.Range("E12").Value2 = 0.6935453139197432
$ws.Range("F12").Value2 = -0.6304082599202043
$ws.Range("G12").Value2 = -0.07823792394983131
$ws.Range("H12").Value2 = 0.6829028338744026
$ws.Range("J12").Value2 = 0.7164471536491567
$ws.Range("K12").Value2 = -0.2857460969514533
$ws.Range("M12").Value2 = 0.2144274027338346
$ws.Range("N12").Value2 = 0.6813215027964721
$ws.Range("O12").Value2 = -0.3284431669736942
$ws.Range("P12").Value2 = 0.6005636968464753
$ws.Range("Q12").Value2 = -0.2324013050705599
$ws.Range("R12").Value2 = 0.8191268319089997
$ws.Range("B13").Value2 = 0.02671256140689233
$ws.Range("C13").Value2 = -0.08960471633066344
$ws.Range("D13").Value2 = -0.03784677341782061
$ws.Range("E13").Value2 = -0.09650071735825808
$ws.Range("F13").Value2 = -0.001045408961199361
$ws.Range("G13").Value2 = 0.813853831783909
$ws.Range("H13").Value2 = 0.48402467427776
$ws.Range("J13").Value2 = 0.4649076623324462
$ws.Range("K13").Value2 = 0.3130512932906556
$ws.Range("L13").Value2 = 0.2144274027338346
$ws.Range("N13").Value2 = 0.4559212458883087
$ws.Range("O13").Value2 = -0.234262498405988
$ws.Range("P13").Value2 = 0.3802176799489076
$ws.Range("Q13").Value2 = 0.3946292024424908
$ws.Range("R13").Value2 = 0.1383314059606824
$ws.Range("B14").Value2 = 0.1245000835920851
$ws.Range("C14").Value2 = -0.1001432587328002
$ws.Range("D14").Value2 = -0.32874929123365
$ws.Range("E14").Value2 = 0.5236914481283271
$ws.Range("F14").Value2 = -0.5212886093024264
$ws.Range("G14").Value2 = 0.1591991380261942
$ws.Range("H14").Value2 = 0.7269516643205035
$ws.Range("J14").Value2 = 0.7604571890318254
$ws.Range("K14").Value2 = -0.06611384108455527
$ws.Range("L14").Value2 = 0.6813215027964721
$ws.Range("M14").Value2 = 0.4559212458883087
$ws.Range("O14").Value2 = -0.3963776518294763
$ws.Range("P14").Value2 = 0.9672529628894712
$ws.Range("Q14").Value2 = 0.01228289385372508
$ws.Range("R14").Value2 = 0.7649321850275498
$ws.Range("B15").Value2 = 0.02613467962342986
$ws.Range("C15").Value2 = 0.01048147561741413
$ws.Range("D15").Value2 = 0.06019978600222384
$ws.Range("E15").Value2 = -0.2827997122318284
$ws.Range("F15").Value2 = 0.270433041677948
$ws.Range("G15").Value2 = -0.1978194138819634
$ws.Range("H15").Value2 = -0.2839099923732194
$ws.Range("J15").Value2 = -0.2806547402725821
$ws.Range("K15").Value2 = 0.03376595019548962
$ws.Range("L15").Value2 = -0.3284431669736942
$ws.Range("M15").Value2 = -0.234262498405988
$ws.Range("N15").Value2 = -0.3963776518294763
$ws.Range("P15").Value2 = -0.3679922888892923
$ws.Range("Q15").Value2 = 0.03262038920980954
$ws.Range("R15").Value2 = -0.3029930305327576
$ws.Range("B16").Value2 = 0.1228681942988932
$ws.Range("C16").Value2 = -0.05947472534720946
$ws.Range("D16").Value2 = -0.2782325710945733
$ws.Range("E16").Value2 = 0.4559128407195984
$ws.Range("F16").Value2 = -0.446859052656183
$ws.Range("G16").Value2 = 0.1469122152520484
$ws.Range("H16").Value2 = 0.6208983952960462
$ws.Range("J16").Value2 = 0.6602733475555668
$ws.Range("K16").Value2 = -0.05246523971468067
$ws.Range("L16").Value2 = 0.6005636968464753
$ws.Range("M16").Value2 = 0.3802176799489076
$ws.Range("N16").Value2 = 0.9672529628894712
$ws.Range("O16").Value2 = -0.3679922888892923
$ws.Range("Q16").Value2 = -0.006289716343228092
$ws.Range("R16").Value2 = 0.7638521715591707
$ws.Range("B17").Value2 = -0.033236964184887
$ws.Range("C17").Value2 = 0.07520151954958959
$ws.Range("D17").Value2 = 0.2684339346906659
$ws.Range("E17").Value2 = -0.6020672565082923
$ws.Range("F17").Value2 = 0.3654188162483002
$ws.Range("G17").Value2 = 0.3630575348992239
$ws.Range("H17").Value2 = 0.04958524081280442
$ws.Range("J17").Value2 = 0.007045408321917977
$ws.Range("K17").Value2 = 0.8799916556406883
$ws.Range("L17").Value2 = -0.2324013050705599
$ws.Range("M17").Value2 = 0.3946292024424908
$ws.Range("N17").Value2 = 0.01228289385372508
$ws.Range("O17").Value2 = 0.03262038920980954
$ws.Range("P17").Value2 = -0.006289716343228092
$ws.Range("R17").Value2 = -0.2590281121190497
$ws.Range("B18").Value2 = 0.001600103591062275
$ws.Range("C18").Value2 = 0.06304668397542296
$ws.Range("D18").Value2 = -0.2432149506306212
$ws.Range("E18").Value2 = 0.5893034284480509
$ws.Range("F18").Value2 = -0.4870047038530387
$ws.Range("G18").Value2 = -0.07265664635526758
$ws.Range("H18").Value2 = 0.527323953863141
$ws.Range("J18").Value2 = 0.5855897097062018
$ws.Range("K18").Value2 = -0.2643290939009635
$ws.Range("L18").Value2 = 0.8191268319089997
$ws.Range("M18").Value2 = 0.1383314059606824
$ws.Range("N18").Value2 = 0.7649321850275498
$ws.Range("O18").Value2 = -0.3029930305327576
$ws.Range("P18").Value2 = 0.7638521715591707
$ws.Range("Q18").Value2 = -0.2590281121190497
